$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the naming-convention values in the last data row (row 14) ---
# Column E (category/run_set label) and Column C (directory) are updated to
# reflect the new naming convention. Set E14 before C14 so new shared
# strings are appended to sharedStrings.xml in the same order as the target.
$ws.Range("E14").Value = "Blueprint"
$ws.Range("C14").Value = "2035_TM152_NGF_NP02_Blueprint_00"

# --- Turn the asana-task cell J14 into a real hyperlink ---
# J14 already holds the asana URL as text; wire it up as a hyperlink
# (this also applies the built-in "Hyperlink" cell style).
$ws.Hyperlinks.Add($ws.Range("J14"), "https://app.asana.com/0/0/1202521542566668/f")

# --- Update the view: scroll right a bit and move the selection to J14 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 1
$ws.Range("J14").Select()

$wb.Save()
